$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.904.94"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.908.33"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.06"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07730"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9802"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.05"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "1.900.69"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.938"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.664"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07026"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.82"
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009458"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.68"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "28.891.42"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.319"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.096"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.53"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.04"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.659"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.54"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.854"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09284"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8673"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.080"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.251"
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.128"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05717"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.163"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02043"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5486"
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.407"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1752"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.863"
$ws.Range("E42").Value = "  +4.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.300"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5175"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.28"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06918"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.085"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002572"
$ws.Range("E48").Value = "  -7.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.776"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.38"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2873"
$ws.Range("E51").Value = "  -4.59%  "
